# Apply updated betting simulation metrics across the three grouped sheets:
# By_Odds_Bin, By_Field_Size, and By_Track.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: By_Odds_Bin
# ---------------------------------------------------------------------------
$wsOdds = $wb.Worksheets.Item("By_Odds_Bin")

# Row 2 : (0, 5]
$wsOdds.Range("B2").Value = 67
$wsOdds.Range("C2").Value = 10.39
$wsOdds.Range("D2").Value = 55.38999999999999
$wsOdds.Range("E2").Value = -45
$wsOdds.Range("F2").Value = 32.8

# Row 3 : (5, 10]
$wsOdds.Range("B3").Value = 32
$wsOdds.Range("C3").Value = 14
$wsOdds.Range("D3").Value = 39
$wsOdds.Range("E3").Value = -25
$wsOdds.Range("F3").Value = 21.9

# Row 4 : (10, 15]
$wsOdds.Range("B4").Value = 1
$wsOdds.Range("C4").Value = -1
$wsOdds.Range("D4").Value = 0
$wsOdds.Range("E4").Value = -1
$wsOdds.Range("F4").Value = 0

# ---------------------------------------------------------------------------
# Sheet 2: By_Field_Size
# ---------------------------------------------------------------------------
$wsField = $wb.Worksheets.Item("By_Field_Size")

# Row 2 : 1-4
$wsField.Range("B2").Value = 13
$wsField.Range("C2").Value = 2.63
$wsField.Range("D2").Value = 11.63
$wsField.Range("E2").Value = -9
$wsField.Range("F2").Value = 30.8

# Row 3 : 5
$wsField.Range("B3").Value = 38
$wsField.Range("C3").Value = 8.68
$wsField.Range("D3").Value = 33.68
$wsField.Range("E3").Value = -25
$wsField.Range("F3").Value = 34.2

# Row 4 : 6
$wsField.Range("B4").Value = 49
$wsField.Range("C4").Value = 12.08
$wsField.Range("D4").Value = 49.08
$wsField.Range("E4").Value = -37
$wsField.Range("F4").Value = 24.5

# The "9-10" field-size bin is split into two separate single-value bins
# ("9" and "10"), and a brand new "11-13" bin is inserted ahead of the
# existing "14+" bin (which shifts from row 9 down to row 10).
# A leading quote forces these numeric-looking labels to stay text labels
# (matching the other bin labels such as "5", "6", "7", "8") instead of
# being auto-converted to numbers.
$wsField.Range("A7").Value = "'9"
$wsField.Range("A8").Value = "'10"

# Insert a new row for the "11-13" bin before the current "14+" row, then
# copy the formatting from the row above it so the label cell keeps the
# same bold/bordered/centered style used throughout the table.
$wsField.Rows.Item(9).Insert()
$wsField.Range("A8").Copy()
$wsField.Range("A9").PasteSpecial(-4122) | Out-Null
$wsField.Application.CutCopyMode = $false

$wsField.Range("A9").Value = "11–13"
$wsField.Range("B9").Value = 0
$wsField.Range("C9").Value = 0
$wsField.Range("D9").Value = 0
$wsField.Range("E9").Value = 0

# ---------------------------------------------------------------------------
# Sheet 3: By_Track
# ---------------------------------------------------------------------------
$wsTrack = $wb.Worksheets.Item("By_Track")

$wsTrack.Range("A2").Value = "HAMILTON"
$wsTrack.Range("B2").Value = 100
$wsTrack.Range("C2").Value = 23.39
$wsTrack.Range("D2").Value = 94.39
$wsTrack.Range("E2").Value = -71
$wsTrack.Range("F2").Value = 29
